# Insert a new data row at row 33 (pushing the existing rows 33-126 down to 34-127)
# and populate it with the new "Early Burlat" cherry price entry for
# Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 33; this shifts rows 33..126 -> 34..127
$ws.Rows.Item(33).EntireRow.Insert()

# Fill in the new row 33 with the new record's values
$ws.Cells.Item(33, 1).Value = 4
$ws.Cells.Item(33, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(33, 3).Value = "Los Lagos"
$ws.Cells.Item(33, 4).Value = 45251
$ws.Cells.Item(33, 5).Value = 10
$ws.Cells.Item(33, 6).Value = "Fruta"
$ws.Cells.Item(33, 7).Value = 100103
$ws.Cells.Item(33, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(33, 9).Value = 100103001
$ws.Cells.Item(33, 10).Value = "Cereza"
$ws.Cells.Item(33, 11).Value = "Early Burlat"
$ws.Cells.Item(33, 12).Value = "Primera"
$ws.Cells.Item(33, 13).Value = 300
$ws.Cells.Item(33, 14).Value = 26000
$ws.Cells.Item(33, 15).Value = 26000
$ws.Cells.Item(33, 16).Value = 26000
$ws.Cells.Item(33, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(33, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(33, 19).Value = 2600
$ws.Cells.Item(33, 20).Value = 10
